# docs/excel-import/samples/Invest_And_Securities_Landscape.xlsx
# doc - improve IMPORT & METAMODEL
#
# Adds two new message-flow rows (TRAD.006 / TRAD.007 - customer
# synchronization between the Core Banking System and the Trading
# Platform / a new Account Microservice), widens the "Step description"
# column so the longer text fits, and leaves the selection on the cell
# below/after the newly entered data (as LibreOffice/Excel would after
# typing the last value and pressing Enter/Tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Message_Flow")

# --- Row 10: TRAD.006 - Core Banking Sysetm -> Trading Platform ---
$ws.Range("A10").Value = "TRAD.006"
$ws.Range("B10").Value = "S.04"
$ws.Range("C10").Value = "Core Banking Sysetm"
$ws.Range("D10").Value = "Trading Platform"
$ws.Range("E10").Value = "Customer synchronization"
$ws.Range("F10").Value = "Send Customer in NRT"
$ws.Range("G10").Value = "Event"
$ws.Range("H10").Value = "NRT"
$ws.Range("I10").Value = "XML"

# --- Row 11: TRAD.007 - Core Banking Sysetm -> Account Microservice ---
$ws.Range("A11").Value = "TRAD.007"
$ws.Range("B11").Value = "S.04"
$ws.Range("C11").Value = "Core Banking Sysetm"
$ws.Range("D11").Value = "Account Microservice"
$ws.Range("E11").Value = "Customer synchronization"
$ws.Range("F11").Value = "Send Customer in batch mode"
$ws.Range("G11").Value = "File"
$ws.Range("H11").Value = "Daily"
$ws.Range("I11").Value = "CSV"

# Widen column F ("Step description") so the newly added, longer text fits.
$ws.Columns.Item(6).ColumnWidth = 62.59

# Leave the cursor where the user would land after entering the new data.
$null = $ws.Range("F12").Select()
